# v1.5 - implementacao de uma rotina para limpar as areas de entrada de informacao
# Adds the new truck-entry (NF 355867 / placa RXQ9H93 / motorista JOAO RAMOS DE
# OLIVEIRA NETO) to "Programacao", splits it across two pallet lots in
# "Planilha", and refreshes the printable "Descarga do Sal" form so it shows
# the newest entry's data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Programacao")
$ws2 = $wb.Worksheets.Item("Planilha")
$ws3 = $wb.Worksheets.Item("Descarga do Sal")

# Text-guard: force text storage (no auto number/date coercion) for values
# that would otherwise be re-interpreted by Excel (pure numeric strings,
# dd/mm/yyyy-looking dates, ...). A leading apostrophe is a no-op for plain
# text, so it is applied to every textual value for consistency.
function T([string]$value) {
    return "'" + $value
}

# ---------------------------------------------------------------------
# Sheet "Programacao": new row 19 - the truck entry itself
# ---------------------------------------------------------------------
$ws1.Range("A19").Value = T("12/02/2025")
$ws1.Range("B19").Value = T("14:46")
$ws1.Range("C19").Value = T("JOAO RAMOS DE OLIVEIRA NETO")
$ws1.Range("D19").Value = T("(84)991067575")
$ws1.Range("E19").Value = T("355867")
$ws1.Range("F19").Value = T("NORSAL")
$ws1.Range("G19").Value = 5400
$ws1.Range("H19").Value = T("RXQ9H93")
$ws1.Range("I19").Value = T("BAU")
$ws1.Range("J19").Value = T("SAL REFINADO 25 KG")
$ws1.Range("K19").Value = T("MONTE SERENO")

# ---------------------------------------------------------------------
# Sheet "Planilha": new rows 30 and 31 - the two pallet-lot splits for
# the NF above (4 pallets + 18 pallets = 22 pallets / 5400 + 24300 = peso)
# ---------------------------------------------------------------------
$ws2.Range("A30").Value = T("ENTRADA")
$ws2.Range("B30").Value = T("12/02/2025")
$ws2.Range("C30").Value = T("RXQ9H93")
$ws2.Range("D30").Value = T("BAU")
$ws2.Range("E30").Value = T("MONTE SERENO")
$ws2.Range("F30").Value = T("SAL REFINADO")
$ws2.Range("G30").Value = T("25 KG")
$ws2.Range("H30").Value = T("NORSAL")
$ws2.Range("I30").Value = T("355867")
$ws2.Range("J30").Value = T("355866")
$ws2.Range("K30").Value = 4
$ws2.Range("L30").Value = T("564303725")
$ws2.Range("M30").Value = T("jan/27")
$ws2.Range("N30").Value = 5400

$ws2.Range("A31").Value = T("ENTRADA")
$ws2.Range("B31").Value = T("12/02/2025")
$ws2.Range("C31").Value = T("RXQ9H93")
$ws2.Range("D31").Value = T("BAU")
$ws2.Range("E31").Value = T("MONTE SERENO")
$ws2.Range("F31").Value = T("SAL REFINADO")
$ws2.Range("G31").Value = T("25 KG")
$ws2.Range("H31").Value = T("NORSAL")
$ws2.Range("I31").Value = T("355869")
$ws2.Range("J31").Value = T("355868")
$ws2.Range("K31").Value = 18
$ws2.Range("L31").Value = T("564302825")
$ws2.Range("M31").Value = T("jan/27")
$ws2.Range("N31").Value = 24300

# ---------------------------------------------------------------------
# Sheet "Descarga do Sal": refresh the form so it reflects the newest
# entry (previously it mirrored the DAO0375 / USINA ALTA MOGIANA entry).
# ---------------------------------------------------------------------
$ws3.Range("K8").Value = T("14:46")
$ws3.Range("D10").Value = T("JOAO RAMOS DE OLIVEIRA NETO")
$ws3.Range("D12").Value = T("(84)991067575")
$ws3.Range("D14").Value = T("BAU")
$ws3.Range("K14").Value = T("RXQ9H93")
$ws3.Range("D16").Value = T("MONTE SERENO")
$ws3.Range("D18").Value = T("NORSAL")
$ws3.Range("M18").Value = T("25 KG")

$ws3.Range("D20").Value = T("355867")
$ws3.Range("K20").Value = T("355866")
$ws3.Range("P20").Value = 5400

$ws3.Range("D22").Value = T("355869")
$ws3.Range("K22").Value = T("355868")
$ws3.Range("P22").Value = 24300

$ws3.Range("D26").Value = T("SAL REFINADO")
$ws3.Range("L26").Value = T("jan/27")

$ws3.Range("D28").Value = T("564303725")
$ws3.Range("H28").Value = T("355867")
$ws3.Range("K28").Value = 5400
$ws3.Range("O28").Value = 4

$ws3.Range("D30").Value = T("564302825")
$ws3.Range("H30").Value = T("355869")
$ws3.Range("K30").Value = 24300
$ws3.Range("O30").Value = 18
